$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update input values
$ws.Range("H4").Value = 10
$ws.Range("H5").Value = 2000

# Update formulas
$ws.Range("C11").Formula = "=IF(H5>500,9.99,3.99)"
$ws.Range("C15").Formula = "=IF((H4*H5)<4000,0,C16)"
$ws.Range("C18").Formula = "=IF((H4*H5)<4000,0,C19)"

# Update selection
$ws.Range("H6").Select()
